# issue #5: stock data from json to db
#
# On the 股票 (stock) sheet, insert a new "category" column between
# property_category and date (always "normal" for this source), and append
# "source_file" (the tmp-file the row was scraped from) and "index" (the
# original per-row id, mirroring column A) columns at the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# --- Insert "category" column before "date" (column I) ---------------------
# Copy column H (property_category) so the freshly inserted column inherits
# the same header/data cell styles used everywhere else on the sheet, then
# overwrite the copied values with the real "category" content.
$ws.Columns.Item(8).Copy()
$ws.Columns.Item(9).Insert()

$ws.Cells.Item(1, 9).Value = "category"
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 9).Value = "normal"
}

# --- Append "source_file" and "index" columns at the end --------------------
$lastCol = $ws.Cells.Item(1, $ws.Columns.Count).End(-4159).Column

$ws.Columns.Item($lastCol).Copy()
$ws.Columns.Item($lastCol + 1).Insert()
$ws.Cells.Item(1, $lastCol + 1).Value = "source_file"
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, $lastCol + 1).Value = "tmpf3421"
}

$lastCol = $ws.Cells.Item(1, $ws.Columns.Count).End(-4159).Column

$ws.Columns.Item($lastCol).Copy()
$ws.Columns.Item($lastCol + 1).Insert()
$ws.Cells.Item(1, $lastCol + 1).Value = "index"
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, $lastCol + 1).Value = $ws.Cells.Item($r, 1).Value2
}
